$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteValues = -4163; used so numeric-looking price strings (e.g. "320.29")
# are written back as literal text, matching the source inlineStr cells,
# instead of being auto-converted to numbers by a plain .Value assignment.
$xlPasteValues = -4163

$ws.Range("D2").Formula = '="41.710.35"'
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial($xlPasteValues)
$ws.Range("E2").Value = '  +0.25%  '

$ws.Range("D3").Formula = '="2.476.44"'
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial($xlPasteValues)
$ws.Range("E3").Value = '  +0.67%  '

$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").Formula = '="320.29"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial($xlPasteValues)
$ws.Range("E5").Value = '  +1.78%  '

$ws.Range("D6").Formula = '="92.34"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial($xlPasteValues)
$ws.Range("E6").Value = '  +0.24%  '

$ws.Range("E7").Value = '  +0.77%  '

$ws.Range("E8").Value = '  +0.05%  '

$ws.Range("D9").Formula = '="0.511"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial($xlPasteValues)
$ws.Range("E9").Value = '  +0.12%  '

$ws.Range("E10").Value = '  +5.44%  '

$ws.Range("D11").Formula = '="33.05"'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial($xlPasteValues)
$ws.Range("E11").Value = '  +2.03%  '

$ws.Range("E12").Value = '  -0.75%  '

$ws.Range("D13").Formula = '="2.857.56"'
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial($xlPasteValues)
$ws.Range("E13").Value = '  +0.67%  '

$ws.Range("E14").Value = '  +0.74%  '

$ws.Range("E15").Value = '  -1.80%  '

$ws.Range("D16").Formula = '="2.486.08"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial($xlPasteValues)
$ws.Range("E16").Value = '  -0.32%  '

$ws.Range("E17").Value = '  +2.23%  '

$ws.Range("D18").Formula = '="41.626.69"'
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial($xlPasteValues)
$ws.Range("E18").Value = '  +0.00%  '

$ws.Range("E19").Value = '  -0.40%  '

$ws.Range("D20").Formula = '="0.0₃0943"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial($xlPasteValues)
$ws.Range("E20").Value = '  -0.05%  '

$ws.Range("D21").Formula = '="70.67"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial($xlPasteValues)
$ws.Range("E21").Value = '  -0.01%  '

$ws.Range("D22").Formula = '="11.27"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial($xlPasteValues)
$ws.Range("E22").Value = '  -0.45%  '

$ws.Range("D23").Formula = '="239.84"'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial($xlPasteValues)
$ws.Range("E23").Value = '  +0.64%  '

$ws.Range("E24").Value = '  +1.55%  '

$ws.Range("E25").Value = '  +2.10%  '

$ws.Range("D26").Formula = '="0.999"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial($xlPasteValues)
$ws.Range("E26").Value = '  -0.05%  '

$ws.Range("D27").Formula = '="25.02"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial($xlPasteValues)
$ws.Range("E27").Value = '  +2.79%  '

$ws.Range("E28").Value = '  -0.64%  '

$ws.Range("E29").Value = '  +0.51%  '

$ws.Range("D30").Formula = '="36.74"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial($xlPasteValues)
$ws.Range("E30").Value = '  +4.65%  '

$ws.Range("D31").Formula = '="157.65"'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial($xlPasteValues)
$ws.Range("E31").Value = '  +1.20%  '

$ws.Range("E32").Value = '  -0.31%  '

$ws.Range("E33").Value = '  +0.00%  '

$ws.Range("D34").Formula = '="0.0764"'
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial($xlPasteValues)
$ws.Range("E34").Value = '  +0.92%  '

$ws.Range("E35").Value = '  +0.02%  '

$ws.Range("D36").Formula = '="17.20"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial($xlPasteValues)
$ws.Range("E36").Value = '  -1.26%  '

$ws.Range("E37").Value = '  +3.00%  '

$ws.Range("E38").Value = '  +1.61%  '

$ws.Range("D39").Formula = '="2.88"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial($xlPasteValues)
$ws.Range("E39").Value = '  -0.04%  '

$ws.Range("E40").Value = '  +1.01%  '

$ws.Range("E41").Value = '  +2.64%  '

$ws.Range("E42").Value = '  -1.76%  '

$ws.Range("D43").Formula = '="1.994.57"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial($xlPasteValues)
$ws.Range("E43").Value = '  +0.93%  '

$ws.Range("E44").Value = '  +0.65%  '

$ws.Range("D45").Formula = '="18.77"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial($xlPasteValues)
$ws.Range("E45").Value = '  -0.33%  '

$ws.Range("E46").Value = '  +2.06%  '

$ws.Range("E47").Value = '  +5.44%  '

$ws.Range("D48").Formula = '="2.755.48"'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial($xlPasteValues)
$ws.Range("E48").Value = '  +2.17%  '

$ws.Range("E49").Value = '  +0.84%  '

$ws.Range("D50").Formula = '="76.01"'
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial($xlPasteValues)
$ws.Range("E50").Value = '  +5.60%  '

$ws.Range("D51").Formula = '="67.51"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial($xlPasteValues)
$ws.Range("E51").Value = '  +1.12%  '
